$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new name "Krishna " into C7 (trailing space preserved)
$ws.Range("C7").Value = "Krishna "

# Clear E7 (previously an empty styled cell; now removed entirely from the sheet)
$ws.Range("E7").ClearFormats()
$ws.Range("E7").ClearContents()

# Update the saved selection to C17
$ws.Range("C17").Select() | Out-Null
